# Updates the "F" column ("想去人数" / interested-count) for the rows
# listed in the commit diff, across all four worksheets.
$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 158  # F3: 157 -> 158
$ws.Cells.Item(6, 6).Value = 315  # F6: 314 -> 315
$ws.Cells.Item(7, 6).Value = 5662  # F7: 5637 -> 5662
$ws.Cells.Item(9, 6).Value = 7634  # F9: 7608 -> 7634
$ws.Cells.Item(10, 6).Value = 302  # F10: 200 -> 302
$ws.Cells.Item(13, 6).Value = 3854  # F13: 3838 -> 3854
$ws.Cells.Item(14, 6).Value = 22  # F14: 21 -> 22
$ws.Cells.Item(15, 6).Value = 26  # F15: 25 -> 26
$ws.Cells.Item(16, 6).Value = 201  # F16: 198 -> 201
$ws.Cells.Item(19, 6).Value = 106  # F19: 105 -> 106
$ws.Cells.Item(21, 6).Value = 605  # F21: 601 -> 605
$ws.Cells.Item(22, 6).Value = 3889  # F22: 3879 -> 3889
$ws.Cells.Item(23, 6).Value = 134  # F23: 133 -> 134
$ws.Cells.Item(25, 6).Value = 5309  # F25: 5286 -> 5309
$ws.Cells.Item(26, 6).Value = 440  # F26: 439 -> 440
$ws.Cells.Item(27, 6).Value = 2104  # F27: 2092 -> 2104
$ws.Cells.Item(29, 6).Value = 353  # F29: 349 -> 353
$ws.Cells.Item(30, 6).Value = 7886  # F30: 7854 -> 7886
$ws.Cells.Item(31, 6).Value = 32  # F31: 31 -> 32
$ws.Cells.Item(33, 6).Value = 2198  # F33: 2194 -> 2198
$ws.Cells.Item(34, 6).Value = 2192  # F34: 2180 -> 2192
$ws.Cells.Item(35, 6).Value = 1334  # F35: 1330 -> 1334
$ws.Cells.Item(36, 6).Value = 1296  # F36: 1293 -> 1296
$ws.Cells.Item(37, 6).Value = 20  # F37: 19 -> 20
$ws.Cells.Item(38, 6).Value = 21  # F38: 20 -> 21
$ws.Cells.Item(39, 6).Value = 267  # F39: 266 -> 267
$ws.Cells.Item(40, 6).Value = 248  # F40: 247 -> 248
$ws.Cells.Item(41, 6).Value = 17  # F41: 15 -> 17
$ws.Cells.Item(42, 6).Value = 1177  # F42: 1176 -> 1177
$ws.Cells.Item(43, 6).Value = 1175  # F43: 1174 -> 1175
$ws.Cells.Item(44, 6).Value = 34  # F44: 32 -> 34
$ws.Cells.Item(45, 6).Value = 1327  # F45: 1325 -> 1327
$ws.Cells.Item(46, 6).Value = 2084  # F46: 2070 -> 2084
$ws.Cells.Item(47, 6).Value = 129  # F47: 127 -> 129
$ws.Cells.Item(48, 6).Value = 221  # F48: 219 -> 221
$ws.Cells.Item(49, 6).Value = 1218  # F49: 1217 -> 1218

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 146  # F4: 145 -> 146

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 571  # F2: 569 -> 571
$ws.Cells.Item(3, 6).Value = 747  # F3: 746 -> 747
$ws.Cells.Item(4, 6).Value = 65  # F4: 64 -> 65

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 158  # F3: 157 -> 158
$ws.Cells.Item(5, 6).Value = 571  # F5: 569 -> 571
$ws.Cells.Item(6, 6).Value = 747  # F6: 746 -> 747
$ws.Cells.Item(7, 6).Value = 315  # F7: 314 -> 315
$ws.Cells.Item(8, 6).Value = 5662  # F8: 5637 -> 5662
$ws.Cells.Item(9, 6).Value = 7635  # F9: 7608 -> 7635
$ws.Cells.Item(10, 6).Value = 303  # F10: 200 -> 303
$ws.Cells.Item(11, 6).Value = 3854  # F11: 3838 -> 3854
$ws.Cells.Item(12, 6).Value = 22  # F12: 21 -> 22
$ws.Cells.Item(13, 6).Value = 26  # F13: 25 -> 26
$ws.Cells.Item(14, 6).Value = 201  # F14: 198 -> 201
$ws.Cells.Item(17, 6).Value = 106  # F17: 105 -> 106
$ws.Cells.Item(19, 6).Value = 146  # F19: 145 -> 146
$ws.Cells.Item(20, 6).Value = 605  # F20: 601 -> 605
$ws.Cells.Item(21, 6).Value = 3889  # F21: 3879 -> 3889
$ws.Cells.Item(23, 6).Value = 134  # F23: 133 -> 134
$ws.Cells.Item(25, 6).Value = 5309  # F25: 5286 -> 5309
$ws.Cells.Item(26, 6).Value = 440  # F26: 439 -> 440
$ws.Cells.Item(27, 6).Value = 2104  # F27: 2092 -> 2104
$ws.Cells.Item(29, 6).Value = 353  # F29: 349 -> 353
$ws.Cells.Item(30, 6).Value = 7886  # F30: 7855 -> 7886
$ws.Cells.Item(31, 6).Value = 32  # F31: 31 -> 32
$ws.Cells.Item(33, 6).Value = 2198  # F33: 2194 -> 2198
$ws.Cells.Item(34, 6).Value = 2192  # F34: 2180 -> 2192
$ws.Cells.Item(35, 6).Value = 1334  # F35: 1330 -> 1334
$ws.Cells.Item(36, 6).Value = 1296  # F36: 1293 -> 1296
$ws.Cells.Item(37, 6).Value = 267  # F37: 266 -> 267
$ws.Cells.Item(38, 6).Value = 248  # F38: 247 -> 248
$ws.Cells.Item(39, 6).Value = 17  # F39: 15 -> 17
$ws.Cells.Item(40, 6).Value = 1177  # F40: 1176 -> 1177
$ws.Cells.Item(41, 6).Value = 1175  # F41: 1174 -> 1175
$ws.Cells.Item(42, 6).Value = 34  # F42: 32 -> 34
$ws.Cells.Item(43, 6).Value = 1327  # F43: 1325 -> 1327
$ws.Cells.Item(44, 6).Value = 2084  # F44: 2070 -> 2084
$ws.Cells.Item(45, 6).Value = 129  # F45: 127 -> 129
$ws.Cells.Item(47, 6).Value = 221  # F47: 219 -> 221
$ws.Cells.Item(49, 6).Value = 1218  # F49: 1217 -> 1218
